{"js": "// Replace the date line and every \"A\u00d7B=C\" answer cell with its updated value.\n// Each old value is unique in the document, so a plain search+replace per\n// pair is sufficient and keeps all run/paragraph formatting untouched.\nconst replacements = [\n  [\"2024-11-19 Tuesday\", \"2024-11-20 Wednesday\"],\n  [\"474\u00d77=3318\", \"945\u00d74=3780\"],\n  [\"297\u00d75=1485\", \"541\u00d74=2164\"],\n  [\"742\u00d72=1484\", \"263\u00d78=2104\"],\n  [\"282\u00d75=1410\", \"958\u00d76=5748\"],\n  [\"473\u00d77=3311\", \"251\u00d79=2259\"],\n  [\"615\u00d74=2460\", \"543\u00d78=4344\"],\n  [\"256\u00d73=768\", \"279\u00d73=837\"],\n  [\"775\u00d72=1550\", \"217\u00d76=1302\"],\n  [\"646\u00d72=1292\", \"498\u00d73=1494\"],\n  [\"970\u00d78=7760\", \"693\u00d73=2079\"],\n  [\"916\u00d74=3664\", \"249\u00d77=1743\"],\n  [\"324\u00d72=648\", \"426\u00d75=2130\"],\n  [\"352\u00d78=2816\", \"990\u00d76=5940\"],\n  [\"775\u00d79=6975\", \"873\u00d75=4365\"],\n  [\"985\u00d78=7880\", \"866\u00d74=3464\"],\n  [\"898\u00d79=8082\", \"577\u00d76=3462\"],\n  [\"304\u00d75=1520\", \"746\u00d78=5968\"],\n  [\"444\u00d74=1776\", \"249\u00d72=498\"],\n  [\"984\u00d74=3936\", \"601\u00d73=1803\"],\n  [\"594\u00d78=4752\", \"570\u00d76=3420\"],\n  [\"648\u00d73=1944\", \"169\u00d76=1014\"],\n  [\"357\u00d77=2499\", \"287\u00d77=2009\"],\n  [\"702\u00d73=2106\", \"991\u00d76=5946\"],\n  [\"244\u00d79=2196\", \"778\u00d72=1556\"],\n  [\"225\u00d75=1125\", \"955\u00d72=1910\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and every \"A\u00d7B=C\" answer cell with its updated value.\n# Each old value is unique in the document, so Find/Replace per pair is\n# sufficient and keeps all run/paragraph formatting untouched.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-11-19 Tuesday\", \"2024-11-20 Wednesday\"),\n    @(\"474\u00d77=3318\", \"945\u00d74=3780\"),\n    @(\"297\u00d75=1485\", \"541\u00d74=2164\"),\n    @(\"742\u00d72=1484\", \"263\u00d78=2104\"),\n    @(\"282\u00d75=1410\", \"958\u00d76=5748\"),\n    @(\"473\u00d77=3311\", \"251\u00d79=2259\"),\n    @(\"615\u00d74=2460\", \"543\u00d78=4344\"),\n    @(\"256\u00d73=768\", \"279\u00d73=837\"),\n    @(\"775\u00d72=1550\", \"217\u00d76=1302\"),\n    @(\"646\u00d72=1292\", \"498\u00d73=1494\"),\n    @(\"970\u00d78=7760\", \"693\u00d73=2079\"),\n    @(\"916\u00d74=3664\", \"249\u00d77=1743\"),\n    @(\"324\u00d72=648\", \"426\u00d75=2130\"),\n    @(\"352\u00d78=2816\", \"990\u00d76=5940\"),\n    @(\"775\u00d79=6975\", \"873\u00d75=4365\"),\n    @(\"985\u00d78=7880\", \"866\u00d74=3464\"),\n    @(\"898\u00d79=8082\", \"577\u00d76=3462\"),\n    @(\"304\u00d75=1520\", \"746\u00d78=5968\"),\n    @(\"444\u00d74=1776\", \"249\u00d72=498\"),\n    @(\"984\u00d74=3936\", \"601\u00d73=1803\"),\n    @(\"594\u00d78=4752\", \"570\u00d76=3420\"),\n    @(\"648\u00d73=1944\", \"169\u00d76=1014\"),\n    @(\"357\u00d77=2499\", \"287\u00d77=2009\"),\n    @(\"702\u00d73=2106\", \"991\u00d76=5946\"),\n    @(\"244\u00d79=2196\", \"778\u00d72=1556\"),\n    @(\"225\u00d75=1125\", \"955\u00d72=1910\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
